$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'332.96"
$ws.Range("E2").Value = "'1.16%"
$ws.Range("G2").Value = "'6"
$ws.Range("D3").Value = "'43.90"
$ws.Range("E3").Value = "'6.91%"
$ws.Range("G3").Value = "'6"
$ws.Range("D4").Value = "'5.769"
$ws.Range("E4").Value = "'2.78%"
$ws.Range("G4").Value = "'6"
$ws.Range("D5").Value = "'0.08348"
$ws.Range("E5").Value = "'2.21%"
$ws.Range("G5").Value = "'6"
$ws.Range("D6").Value = "'8.814"
$ws.Range("E6").Value = "'0.86%"
$ws.Range("G6").Value = "'6"
$ws.Range("D7").Value = "'1.975"
$ws.Range("E7").Value = "'-3.32%"
$ws.Range("G7").Value = "'6"
$ws.Range("D8").Value = "'2.900"
$ws.Range("E8").Value = "'-2.90%"
$ws.Range("G8").Value = "'6"
$ws.Range("D9").Value = "'0.9337"
$ws.Range("E9").Value = "'1.69%"
$ws.Range("G9").Value = "'6"
$ws.Range("D10").Value = "'0.1241"
$ws.Range("E10").Value = "'-1.74%"
$ws.Range("G10").Value = "'6"
$ws.Range("D11").Value = "'0.1947"
$ws.Range("E11").Value = "'-0.65%"
$ws.Range("G11").Value = "'6"
$ws.Range("D12").Value = "'0.09639"
$ws.Range("E12").Value = "'2.57%"
$ws.Range("G12").Value = "'6"
$ws.Range("D13").Value = "'0.03942"
$ws.Range("E13").Value = "'5.54%"
$ws.Range("G13").Value = "'6"
$ws.Range("D14").Value = "'0.1067"
$ws.Range("E14").Value = "'0.85%"
$ws.Range("G14").Value = "'6"
$ws.Range("D15").Value = "'0.001301"
$ws.Range("E15").Value = "'-0.56%"
$ws.Range("G15").Value = "'6"
$ws.Range("D16").Value = "'0.005928"
$ws.Range("E16").Value = "'-3.94%"
$ws.Range("G16").Value = "'6"
$ws.Range("D17").Value = "'3.504"
$ws.Range("E17").Value = "'2.02%"
$ws.Range("G17").Value = "'6"
$ws.Range("D18").Value = "'4.503"
$ws.Range("E18").Value = "'-0.67%"
$ws.Range("G18").Value = "'6"
$ws.Range("G19").Value = "'6"
$ws.Range("D20").Value = "'9.035"
$ws.Range("E20").Value = "'8.94%"
$ws.Range("G20").Value = "'6"
$ws.Range("D21").Value = "'0.1371"
$ws.Range("E21").Value = "'-1.59%"
$ws.Range("G21").Value = "'6"
$ws.Range("E22").Value = "'7.54%"
$ws.Range("G22").Value = "'6"
$ws.Range("D23").Value = "'0.04412"
$ws.Range("E23").Value = "'-0.06%"
$ws.Range("G23").Value = "'6"
$ws.Range("D24").Value = "'0.001259"
$ws.Range("E24").Value = "'-0.38%"
$ws.Range("G24").Value = "'6"
$ws.Range("D25").Value = "'0.004365"
$ws.Range("E25").Value = "'1.93%"
$ws.Range("G25").Value = "'6"
$ws.Range("G26").Value = "'6"
$ws.Range("D27").Value = "'0.0003993"
$ws.Range("E27").Value = "'0.00%"
$ws.Range("G27").Value = "'6"
$ws.Range("G28").Value = "'6"
$ws.Range("G29").Value = "'6"
$ws.Range("G30").Value = "'6"
$ws.Range("G31").Value = "'6"
$ws.Range("G32").Value = "'6"
$ws.Range("G33").Value = "'6"
$ws.Range("G34").Value = "'6"
$ws.Range("G35").Value = "'6"
$ws.Range("G36").Value = "'6"
$ws.Range("G37").Value = "'6"
$ws.Range("G38").Value = "'6"
$ws.Range("D39").Value = "'0.02840"
$ws.Range("E39").Value = "'3.33%"
$ws.Range("G39").Value = "'6"
$ws.Range("D40").Value = "'0.05742"
$ws.Range("E40").Value = "'6.27%"
$ws.Range("G40").Value = "'6"
$ws.Range("D41").Value = "'0.007920"
$ws.Range("E41").Value = "'3.45%"
$ws.Range("G41").Value = "'6"
$ws.Range("D42").Value = "'0.1427"
$ws.Range("E42").Value = "'0.99%"
$ws.Range("G42").Value = "'6"
$ws.Range("D43").Value = "'0.009077"
$ws.Range("E43").Value = "'0.62%"
$ws.Range("G43").Value = "'6"
$ws.Range("D44").Value = "'0.002102"
$ws.Range("E44").Value = "'-1.51%"
$ws.Range("G44").Value = "'6"
$ws.Range("D45").Value = "'0.01019"
$ws.Range("E45").Value = "'-9.85%"
$ws.Range("G45").Value = "'6"
$ws.Range("D46").Value = "'0.00007291"
$ws.Range("E46").Value = "'5.74%"
$ws.Range("G46").Value = "'6"
$ws.Range("E47").Value = "'-0.12%"
$ws.Range("G47").Value = "'6"
$ws.Range("D48").Value = "'0.003196"
$ws.Range("E48").Value = "'-11.00%"
$ws.Range("G48").Value = "'6"
$ws.Range("E49").Value = "'-0.19%"
$ws.Range("G49").Value = "'6"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.12%"
$ws.Range("G50").Value = "'6"
$ws.Range("E51").Value = "'-0.12%"
$ws.Range("G51").Value = "'6"
